$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'20.561.98"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.59%  '
$ws.Range('D3').Value = "'1.469.22"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E4').Value = '  +0.34%  '
$ws.Range('D5').Value = "'0.9596"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.32%  '
$ws.Range('D6').Value = "'276.64"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.69%  '
$ws.Range('D7').Value = "'0.3549"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.50%  '
$ws.Range('D8').Value = "'0.3056"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.72%  '
$ws.Range('D9').Value = "'1.081"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.33%  '
$ws.Range('D10').Value = "'39.25"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.27%  '
$ws.Range('D11').Value = "'0.06611"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.01%  '
$ws.Range('D13').Value = "'5.444"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.48%  '
$ws.Range('D14').Value = "'18.02"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.63%  '
$ws.Range('D15').Value = "'6.157"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.33%  '
$ws.Range('D16').Value = "'0.9604"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.36%  '
$ws.Range('D17').Value = "'0.00001015"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.01%  '
$ws.Range('D18').Value = "'1.472.48"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.14%  '
$ws.Range('D19').Value = "'0.05936"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.71%  '
$ws.Range('D20').Value = "'68.73"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.93%  '
$ws.Range('D21').Value = "'5.464"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.59%  '
$ws.Range('D22').Value = "'14.41"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.87%  '
$ws.Range('D23').Value = "'11.16"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.31%  '
$ws.Range('D24').Value = "'2.268"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.17%  '
$ws.Range('D25').Value = "'20.581.52"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.53%  '
$ws.Range('D26').Value = "'145.56"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.40%  '
$ws.Range('D27').Value = "'2.079"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.18%  '
$ws.Range('D28').Value = "'17.05"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.03%  '
$ws.Range('D29').Value = "'1.630.12"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.26%  '
$ws.Range('D30').Value = "'113.99"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.53%  '
$ws.Range('D31').Value = "'3.966"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.29%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value = "'0.07943"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.76%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = "'4.899"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.00%  '
$ws.Range('D34').Value = "'0.7891"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.26%  '
$ws.Range('E35').Value = '  +7.86%  '
$ws.Range('D36').Value = "'1.449"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.18%  '
$ws.Range('D37').Value = "'0.05671"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.14%  '
$ws.Range('D38').Value = "'4.696"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.56%  '
$ws.Range('D39').Value = "'0.9606"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.46%  '
$ws.Range('D40').Value = "'0.02027"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.38%  '
$ws.Range('D41').Value = "'10.23"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.12%  '
$ws.Range('D42').Value = "'0.1840"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.40%  '
$ws.Range('D43').Value = "'7.282"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.11%  '
$ws.Range('E44').Value = '  +0.96%  '
$ws.Range('D45').Value = "'0.5203"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.57%  '
$ws.Range('D46').Value = "'11.97"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.09%  '
$ws.Range('D47').Value = "'119.82"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.09%  '
$ws.Range('D48').Value = "'0.5148"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.37%  '
$ws.Range('D50').Value = "'0.06418"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.92%  '
$ws.Range('D51').Value = "'0.9884"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.05%  '
